$wb = $excel.ActiveWorkbook

# --- Sheet "List": drop the "id" column (A) and the "account" column (last) ---
$ws1 = $wb.Worksheets.Item("List")
$ws1.Columns("A:A").Delete()
$ws1.Columns("D:E").Delete()

# --- Sheet "Search": insert a new "id" row before the existing detail rows,
#     and turn the old "account" row (now shifted down) into a "user" row ---
$ws2 = $wb.Worksheets.Item("Search")
$ws2.Rows("4:4").Insert()

$ws2.Range("A4").Value = '${msg.getProperty(''savedSearch_id'')}'
$ws2.Range("B4").Value = '${id}'

$ws2.Range("A8").Value = '${msg.getProperty(''savedSearch_user'')}'
$ws2.Range("B8").Value = '${user}'
